$d = $word.ActiveDocument

# 1. Remove the stray _GoBack bookmark that currently sits at the very
#    start of the title paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Insert a brand-new paragraph right after the "AS" (added subclass
#    (restriction)) paragraph, containing the "CS" / "changed subclass
#    (restriction)" abbreviation entry, with the _GoBack bookmark placed
#    at the very end of that new paragraph (after its text, before the
#    paragraph mark) -- matching where it now belongs per the diff.
$asParagraph = $d.Paragraphs.Item(11)
$asParagraph.Range.InsertParagraphAfter()

$csParagraph = $d.Paragraphs.Item(12)
$newParagraphXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>CS</w:t></w:r><w:r><w:tab/><w:t>changed subclass (restriction)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$csParagraph.Range.InsertXML($newParagraphXml) | Out-Null
